# Add a new data row (row 4) to the "Artfynd" sheet, mirroring the
# structure of the existing rows (2 and 3): a new species observation
# record appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric fields -------------------------------------------------
$ws.Range("A4").Value = 131204756
$ws.Range("B4").Value = 57881
$ws.Range("E4").Value = 100049
$ws.Range("Q4").Value = 604965
$ws.Range("R4").Value = 6546211
$ws.Range("S4").Value = 50

# --- Plain text fields (not number/date-like, safe to assign directly) ---
$ws.Range("D4").Value = "NT"
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "lockläte, övriga läten"
$ws.Range("P4").Value = "Ällmora SV om , Srm"
$ws.Range("T4").Value = "Södermanland"
$ws.Range("U4").Value = "Flen"
$ws.Range("V4").Value = "Södermanland"
$ws.Range("W4").Value = "Helgesta"
$ws.Range("AW4").Value = "Rolf Olsson"
$ws.Range("AX4").Value = "Rolf Olsson"

# --- Text fields that look like numbers/dates: force text formatting so
#     they are stored as text (matching the source data) instead of being
#     auto-coerced into a number / date serial. ---------------------------
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"

$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2026-02-17"

$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2026-02-17"

# --- Boolean fields --------------------------------------------------------
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
